# Update odds-base data: swap the contents of columns B..AC between the
# listed adjacent row pairs (the "id" values in column A stay put — only
# the match/record data that had been misaligned between the two rows is
# exchanged back into the correct row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 29  # column AC

$rowPairs = @(
    @(16, 17),
    @(20, 21),
    @(56, 57),
    @(83, 84),
    @(90, 91),
    @(110, 111),
    @(113, 114)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}
